$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A8").Value = 46042
$ws.Range("D8").Value = 156.52000000000001
$ws.Range("E8").Value = 150.82
$ws.Range("F8").Value = 160.82
$ws.Range("G8").Value = 150.84
$ws.Range("A9").Value = 46042
$ws.Range("D9").Value = 156.52000000000001
$ws.Range("E9").Value = 150.82
$ws.Range("F9").Value = 160.82
$ws.Range("G9").Value = 150.84
$ws.Range("A10").Value = 46042
$ws.Range("D10").Value = 157.54
$ws.Range("E10").Value = 152.82
$ws.Range("F10").Value = 162.82
$ws.Range("G10").Value = 153.22
$ws.Range("A11").Value = 46039
$ws.Range("D11").Value = 156.03
$ws.Range("E11").Value = 150.84
$ws.Range("F11").Value = 160.84
$ws.Range("G11").Value = 150.85
$ws.Range("A12").Value = 46039
$ws.Range("D12").Value = 156.03
$ws.Range("E12").Value = 150.84
$ws.Range("F12").Value = 160.84
$ws.Range("G12").Value = 150.85
$ws.Range("A13").Value = 46039
$ws.Range("D13").Value = 157.01
$ws.Range("E13").Value = 153.01
$ws.Range("F13").Value = 163.01
$ws.Range("G13").Value = 153.41
$ws.Range("A17").Value = 46042
$ws.Range("D17").Value = 160.99
$ws.Range("E17").Value = 155.37
$ws.Range("F17").Value = 165.37
$ws.Range("A18").Value = 46039
$ws.Range("D18").Value = 160.51
$ws.Range("E18").Value = 155.6
$ws.Range("F18").Value = 165.6
$ws.Range("A22").Value = 46042
$ws.Range("D22").Value = 157.71
$ws.Range("E22").Value = 152.51
$ws.Range("F22").Value = 162.11000000000001
$ws.Range("G22").Value = 153.58000000000001
$ws.Range("A23").Value = 46042
$ws.Range("D23").Value = 162.30000000000001
$ws.Range("E23").Value = 158.15
$ws.Range("F23").Value = 168.15
$ws.Range("A24").Value = 46042
$ws.Range("D24").Value = 162.44
$ws.Range("E24").Value = 158.85
$ws.Range("F24").Value = 168.85
$ws.Range("A25").Value = 46042
$ws.Range("D25").Value = 162.43
$ws.Range("E25").Value = 158.38
$ws.Range("F25").Value = 168.38
$ws.Range("G25").Value = 158.51
$ws.Range("A26").Value = 46042
$ws.Range("D26").Value = 162.03
$ws.Range("E26").Value = 160.01
$ws.Range("F26").Value = 170.01
$ws.Range("A27").Value = 46039
$ws.Range("D27").Value = 157.11000000000001
$ws.Range("E27").Value = 152.63
$ws.Range("F27").Value = 162.22999999999999
$ws.Range("G27").Value = 153.69999999999999
$ws.Range("A28").Value = 46039
$ws.Range("D28").Value = 161.78
$ws.Range("E28").Value = 158.35
$ws.Range("F28").Value = 168.35
$ws.Range("A29").Value = 46039
$ws.Range("D29").Value = 161.91999999999999
$ws.Range("E29").Value = 159.02000000000001
$ws.Range("F29").Value = 169.02
$ws.Range("A30").Value = 46039
$ws.Range("D30").Value = 161.9
$ws.Range("E30").Value = 158.56
$ws.Range("F30").Value = 168.56
$ws.Range("G30").Value = 158.68
$ws.Range("A31").Value = 46039
$ws.Range("D31").Value = 161.51
$ws.Range("E31").Value = 160.18
$ws.Range("F31").Value = 170.18
$ws.Range("A35").Value = 46042
$ws.Range("D35").Value = 156.41999999999999
$ws.Range("E35").Value = 149.81
$ws.Range("F35").Value = 158.81
$ws.Range("A36").Value = 46039
$ws.Range("D36").Value = 155.88999999999999
$ws.Range("E36").Value = 149.99
$ws.Range("F36").Value = 158.99
$ws.Range("A40").Value = 46042
$ws.Range("D40").Value = 161.91999999999999
$ws.Range("E40").Value = 158.16
$ws.Range("F40").Value = 168.16
$ws.Range("A41").Value = 46042
$ws.Range("D41").Value = 161.63999999999999
$ws.Range("E41").Value = 158.58000000000001
$ws.Range("F41").Value = 168.58
$ws.Range("A42").Value = 46039
$ws.Range("D42").Value = 161.37
$ws.Range("E42").Value = 158.36000000000001
$ws.Range("F42").Value = 168.36
$ws.Range("A43").Value = 46039
$ws.Range("D43").Value = 161.09
$ws.Range("E43").Value = 158.78
$ws.Range("F43").Value = 168.78
$ws.Range("A47").Value = 46042
$ws.Range("D47").Value = 155.81
$ws.Range("E47").Value = 151.04
$ws.Range("F47").Value = 161.04
$ws.Range("A48").Value = 46042
$ws.Range("D48").Value = 155.41
$ws.Range("E48").Value = 150.96
$ws.Range("F48").Value = 160.96
$ws.Range("A49").Value = 46039
$ws.Range("D49").Value = 154.4
$ws.Range("E49").Value = 150.66999999999999
$ws.Range("F49").Value = 160.66999999999999
$ws.Range("A50").Value = 46039
$ws.Range("D50").Value = 154.02000000000001
$ws.Range("E50").Value = 150.6
$ws.Range("F50").Value = 160.6
$ws.Range("A54").Value = 46042
$ws.Range("D54").Value = 171
$ws.Range("E54").Value = 165.89
$ws.Range("F54").Value = 175.89
$ws.Range("A55").Value = 46042
$ws.Range("D55").Value = 163.68
$ws.Range("E55").Value = 164.08
$ws.Range("F55").Value = 174.08
$ws.Range("A56").Value = 46042
$ws.Range("D56").Value = 160.47
$ws.Range("A57").Value = 46042
$ws.Range("D57").Value = 160.91
$ws.Range("E57").Value = 158.5
$ws.Range("A58").Value = 46042
$ws.Range("D58").Value = 156.66999999999999
$ws.Range("E58").Value = 154.4
$ws.Range("F58").Value = 164.4
$ws.Range("A59").Value = 46042
$ws.Range("D59").Value = 163.31
$ws.Range("E59").Value = 164.01
$ws.Range("A60").Value = 46039
$ws.Range("D60").Value = 170.5
$ws.Range("E60").Value = 165.91
$ws.Range("F60").Value = 175.91
$ws.Range("A61").Value = 46039
$ws.Range("D61").Value = 163.09
$ws.Range("E61").Value = 164.2
$ws.Range("F61").Value = 174.2
$ws.Range("A62").Value = 46039
$ws.Range("D62").Value = 159.83000000000001
$ws.Range("A63").Value = 46039
$ws.Range("D63").Value = 160.32
$ws.Range("E63").Value = 158.62
$ws.Range("A64").Value = 46039
$ws.Range("D64").Value = 156.09
$ws.Range("E64").Value = 154.52000000000001
$ws.Range("F64").Value = 164.52
$ws.Range("A65").Value = 46039
$ws.Range("D65").Value = 162.84
$ws.Range("E65").Value = 164.11
